# Reorder the "Recorded By" (column G) values so that any "System"/"system"
# token(s) are moved to the front of the comma-separated list, while all
# other tokens (e.g. email addresses) keep their original relative order
# and are shifted to the end.
#
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
#          "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $text = [string]$value

    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text -split ","
    $sysParts = @()
    $otherParts = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed -eq "System" -or $trimmed -eq "system") {
            $sysParts += $trimmed
        } else {
            $otherParts += $trimmed
        }
    }

    if ($sysParts.Count -eq 0) {
        continue
    }

    $newParts = $sysParts + $otherParts
    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $text) {
        $cell.Value = $newValue
    }
}
